$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the "Predicted Date of Failure" (G) and "Comment" (H) columns ---
$ws.Columns("G:H").Delete()

# --- New data table (rows 2..15), columns A-F ---
# A: Coach Number, B: Axle Number, C: Wheel ID, D: Defect Size (mm),
# E: Reference Date (kept as the existing "2017-04-07" text), F: Number of Days before Failure
$data = @(
    @(12464, 4, 8, 28,      0),
    @(10805, 4, 8, 30.1317, 256),
    @(10805, 4, 7, 30.064,  259),
    @(10805, 3, 6, 30.2994, 247),
    @(10805, 2, 4, 32.9504, 107),
    @(10805, 2, 3, -4.2395, 0),
    @(10805, 1, 2, 32.1448, 150),
    @(10805, 1, 1, 33.4024, 84),
    @(12464, 1, 1, 31.2552, 197),
    @(12464, 2, 3, 30.1516, 255),
    @(12464, 2, 4, 6.6626,  130),
    @(12464, 3, 5, 30.657,  228),
    @(12464, 3, 6, 30.2795, 248),
    @(12464, 4, 7, 30.2506, 249)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $startRow + $i
    $row = $data[$i]

    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    # Copy the existing text cell so the "2017-04-07" stays a shared string,
    # not an auto-converted date serial number.
    $ws.Range("E2").Copy($ws.Cells.Item($r, 5))
    $ws.Cells.Item($r, 6).Value = $row[4]
}
